$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2025-10-27"

# Copy header row (values + style) from the previous last sheet so the bold/border/
# alignment header style (style index 1) is reused instead of minting a new one.
$lastSheet.Range("A1:D1").Copy($ws.Range("A1:D1"))

$data = @(
  @('1', '異世界おじさん', '殆ど死んでいる(著者)', '第72話'),
  @('2', '新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる', '漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea', '第72話 前編'),
  @('3', '転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～', 'zunta(作画) はらわたさいぞう(原作)', '第33話：誰にも負けない完璧③'),
  @('4', '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―', '光永康則', '第７０話「突貫停止」③'),
  @('5', 'いとこのこ', 'いぬちく(著者)', '第40話'),
  @('6', '生徒会にも穴はある！', 'むちまろ', '第138話	虎丸の世界（後編）'),
  @('7', '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。', 'マツモトケンゴ', '第６６話　登山の戦いが始まった（２）'),
  @('8', 'クラスで２番目に可愛い女の子と友だちになった', '尾野凛(漫画) たかた(原作) 日向あずり(キャラクター原案)', '第36話②'),
  @('9', '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが', '絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)', '第5話 後編'),
  @('10', '男女比1：5の世界でも普通に生きられると思った？　～激重感情な彼女たちが無自覚男子に翻弄されたら～', '三藤 孝太郎(原作) 桃季憂(漫画) jimmy(キャラクター原案)', '第12話-2'),
  @('11', '実は俺、最強でした？', '原作：澄守 彩 漫画：高橋 愛', '第130話　ライアスの苦悩・前編'),
  @('12', '黒幕一家に転生したけど原作無視して独立する', '空野進 赤村晃人 笠間三四郎 るろお', '第５話　ゲス子爵を成敗して独立する（１）'),
  @('13', '異世界のんびり農家', '剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)', '第309話'),
  @('14', '淫獄団地', '搾精研究所(原作) 丈山雄為(漫画)', '第51話（後編）'),
  @('15', '不純な彼女達は懺悔しない', 'ポロロッカ(著者)', '第33話'),
  @('16', '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～', '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり', '第５４話　勇者を護る器用貧乏（１）'),
  @('17', 'リビルドワールド', '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)', '第74話③'),
  @('18', '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！', '島知宏 音速炒飯 有都あらゆる', '番外編②　マリーの日記ですわ！'),
  @('19', 'バキ外伝 烈海王は異世界転生しても一向にかまわんッッ', '板垣恵介 猪原賽 陸井栄史', '第82話　蛇妖女(メデューサ)'),
  @('20', '独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～', '漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき', '第35話 独身貴族はバーでハイボールを作る（2）'),
  @('21', '「ククク……。奴は四天王の中でも最弱」と解雇された俺、なぜか勇者と聖女の師匠になる', '漫画：芳橋アツシ 原作：延野正行 キャラクター原案：坂野杏梨', '第42話　奴は座っている（後編）'),
  @('22', 'よくわからないけれど異世界に転生していたようです', '内々けやき あし カオミン', '第141話 よくわからないけれど超理論が生み出されたようです（２）'),
  @('23', '異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～', 'モリタ Ｕ４ nima', '第14話（１）　春とおぼっちゃまとピクニックランチ（１）'),
  @('24', '美人女上司滝沢さん', 'やんBARU(著者)', '第206話'),
  @('25', 'くらいあの子としたいこと', '碇マナツ(著者)', '特別編⑳'),
  @('26', '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~', '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)', '第83話(前編)その3'),
  @('27', '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信', '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)', '第4話前編'),
  @('28', '衛宮さんちの今日のごはん', 'TAa(漫画) 只野まこと(料理監修) ＴＹＰＥ－ＭＯＯＮ(原作)', '第76話'),
  @('29', '物語の黒幕に転生して', '瀬川はじめ(漫画) 結城涼(原作) なかむら(キャラクター原案)', '第35話'),
  @('30', 'バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～', '板垣恵介 林たかあき', '第56話 暑い夏、再び'),
  @('31', '姫様“拷問”の時間です', '原作:春原ロビンソン　漫画:ひらけい', '拷問154'),
  @('32', 'アザミヤコを好きになる', 'ユニティコング(原作) ツノニガウ(作画)', '単行本2巻豪華応援イラスト！！'),
  @('33', '老後に備えて異世界で８万枚の金貨を貯めます', 'FUNA 東西 モトエ恵介', '第123話　襲撃［その２］'),
  @('34', 'インフィニット・デンドログラム', '今井神 原作：海道左近 キャラクター原案：タイキ', '第74話'),
  @('35', '十年目、帰還を諦めた転移者はいまさら主人公になる', '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう', '第２０話②'),
  @('36', 'Lv２からチートだった元勇者候補のまったり異世界ライフ', '糸町秋音（漫画） 鬼ノ城ミヤ（原作） 片桐（キャラクター原案）', '第62話　双子の成長'),
  @('37', '理想のヒモ生活', '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)', '第88話　その4'),
  @('38', '解雇された暗黒兵士(30代)のスローなセカンドライフ', '岡沢六十四 るれくちぇ sage・ジョー', '第74話(前編) イダ乱入'),
  @('39', 'ふかふかダンジョン攻略記～俺の異世界転生冒険譚～', 'KAKERU', '第69話「光の剣」（後編）'),
  @('40', '俺は星間国家の悪徳領主！', '灘島かい（漫画） 三嶋与夢（原作） 高峰ナダレ（キャラクター原案）', '第41話　究極にして至高（後編）'),
  @('41', '彼女にしたい女子一位、の隣で見つけたあまりちゃん', '寝巻ネルゾ(漫画) 裕時悠示(原作) たん旦(キャラクター原案)', '第6話①「ダブルデート？」'),
  @('42', '魔王になったので、ダンジョン造って人外娘とほのぼのする', '遠野ノオト(作画) 流優(原作) だぶ竜(キャラクター原案)', '第84話後半'),
  @('43', '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版', '漫画/すたひろ 原作/Y.A', 'chapter72【38話①】'),
  @('44', '幼女戦記', '東條チカ(漫画) カルロ・ゼン(原作) 篠月しのぶ(キャラクター原案)', '第百九章：ドードーバード航空戦Ⅳ'),
  @('45', '塔の管理をしてみよう', '盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)', '第94話前編'),
  @('46', '最弱貴族に転生したので悪役たちを集めてみた', '空野進 sorani ファルまろ', '第13話　最弱貴族、悪役令嬢に賭けを挑む（３）'),
  @('47', '女子高生の無駄づかい', 'ビーノ(著者)', '第137話　せいとかい'),
  @('48', '黄金の経験値', '原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)', '第19話（前編）'),
  @('49', 'オタクに優しいギャルはいない!?', 'のりしろちゃん 魚住さかな', '【特別話】オタクに優しいギャルはいない!?_学生証風ICカードステッカーPR'),
  @('50', '異世界のすみっこで快適ものづくり生活 ～女神さまのくれた工房はちょっとやりすぎ性能だった～', '西山アラタ(漫画) 長田信織(原作) 東上文(キャラクター原案)', 'EP.20②')
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = [int]$row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}
